# Apache POI bug 59983: add a new summary row (row 6) below the existing
# "Total" row, summing the three rows above it (rows 3-5), so the test
# fixture exercises shared-formula handling when rows are later shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C6 is a standalone SUM formula.
$ws.Range("C6").Formula = "=SUM(C3:C5)"

# D6 and E6 share one relative formula pattern, entered in a single
# operation so Excel records them as a shared-formula group (mirrors the
# existing D5:E5 shared formula already in the sheet).
$ws.Range("D6:E6").Formula = "=SUM(D3:D5)"

# Move/restore the active selection to C7, just below the new row.
$ws.Range("C7").Select()
